# Word COM-interop script reproducing the tracked edit to
# "Document de presentation.docx".
#
# Summary of the change:
#   1. Paragraph "Le projet realise sera un jeu de role (RPG) en 3D."
#      gets a left tab stop at 6105 twips and a new sentence appended
#      as a separate run.
#   2. Paragraph "Le dit jeu de role aura un theme bien precis..." is
#      reworded/extended and split into two paragraphs right after
#      "...augmenter son niveau. ".
#   3. The _GoBack bookmark moves from the "Cas d'usage" placeholder
#      paragraph to the new split point created in step 2.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Le projet realise..." paragraph: add a tab stop + append a run.
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(3)
$p1.Range.ParagraphFormat.TabStops.Add(6105 / 20)

$p1 = $d.Paragraphs(3)
$end1 = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$end1.Collapse(0)
$end1.InsertAfter(" L’objectif recherché en entamant le développement de ce dernier est de m’entraîner à utiliser les principes de polymorphisme et d’encapsulation des données dans un cadre de programmation des plus concrets possibles.")

# ---------------------------------------------------------------------
# 2) "Le dit jeu de role..." paragraph: reword + split.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "nouveaux items, gagner de l’expérience afin d’augmenter son niveau, ou encore régénérer",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "nouveaux items, afin d’être mieux équipé pour affronter les prédateurs fauniques ainsi que les tribus locales hostiles; ceci sera nécessaire afin de gagner de l’expérience et enfin d’augmenter son niveau. ^pu encore régénérer",
    2)

# trailing space added after "(HP)."
$d.Content.Find.Execute(
    "régénérer ses points de vie (HP).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "régénérer ses points de vie (HP). ",
    2)

# ---------------------------------------------------------------------
# 3) Move the _GoBack bookmark to the new split point.
# ---------------------------------------------------------------------
$old = $d.Bookmarks("_GoBack")
$old.Delete()

$d.Content.Find.Execute("augmenter son niveau. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$found = $d.Content.Find.Parent
$markerEnd = $found.End
$markerRange = $d.Range($markerEnd, $markerEnd)
$d.Bookmarks.Add("_GoBack", $markerRange)
